# Deploy website Tue Nov 15 16:57:06 PST 2022
#
# 1) Slide 17 ("The Web: HTML"): the two runs " partial section of the "
#    and "CS88 Website:" get merged into a single run (keeping the
#    formatting/rPr of the first run).
# 2) Slide 19 ("Prolog Example"): the wide diagram picture is nudged to a
#    new position.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 17: merge " partial section of the " + "CS88 Website:" into one run
# ---------------------------------------------------------------------
$s17   = $p.Slides.Item(17)
$shp17 = $s17.Shapes.Item(2)            # "Content Placeholder 2"
$tr17  = $shp17.TextFrame.TextRange

$full17    = $tr17.Text
$needle    = " partial section of the "
$startIdx  = $full17.IndexOf($needle)
$mergedStr = " partial section of the CS88 Website:"

# Re-assigning the identical text across both original runs collapses
# them into a single run, carrying over the formatting of the run that
# starts at $startIdx (i.e. the " partial section of the " run).
$mergeRange = $tr17.Characters($startIdx + 1, $mergedStr.Length)
$mergeRange.Text = $mergeRange.Text

# ---------------------------------------------------------------------
# Slide 19: reposition the wide diagram picture (Google Shape;156;p25)
# ---------------------------------------------------------------------
$s19   = $p.Slides.Item(19)
$shp19 = $s19.Shapes.Item(4)            # "Google Shape;156;p25"

$emuPerPt   = 914400.0 / 72.0
$targetXEmu = 3809627
$targetYEmu = 1234402

# Shape.Left/.Top are single-precision (Single) in the COM object model,
# so converting EMU -> points and back can land one EMU short after
# float32 rounding; nudge by a tiny epsilon so the stored value lands on
# the intended EMU after PowerPoint's internal float32 round-trip.
$shp19.Left = $targetXEmu / $emuPerPt
$shp19.Top  = ($targetYEmu / $emuPerPt) + 0.00003
